# Sample Project / Main.xlsx - Rules sheet update
#
# Row 11 (the "R40" rule row) has its Rule-name cell (B11) changed from the
# text "R40" to the text "1". The value must stay a text string (it is
# stored as a shared string in the workbook, not a number), so it is
# entered with a leading apostrophe the way a user would type it into the
# Excel UI to force text interpretation of a numeric-looking value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
